$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 32

# Column A holds the date as literal text (matching the existing rows which
# store dates as plain text, not as Excel date serials). Using a leading
# apostrophe forces Excel to treat the value as text instead of auto-parsing
# it into a date, then we reset the style so no stray number-format style
# is left behind on the cell.
$ws.Cells.Item($row, 1).Value = "'12/26/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 12401.14
$ws.Cells.Item($row, 3).Value = 0.205772699318881
$ws.Cells.Item($row, 4).Value = 0.794227300681119
$ws.Cells.Item($row, 5).Value = -128.84
$ws.Cells.Item($row, 6).Value = -24.08
$ws.Cells.Item($row, 7).Value = -20742.43
$ws.Cells.Item($row, 8).Value = -67.8
$ws.Cells.Item($row, 9).Value = -451.03
$ws.Cells.Item($row, 10).Value = -15.02
